$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 169, pushing the existing
# rows 169-196 down to 172-199 (this also extends the used range to
# A1:T199, matching the updated <dimension> in the diff).
$ws.Rows("169:171").Insert()

# Row 169: new weekly price entry (Especial)
$ws.Range("A169").Value = 6
$ws.Range("B169").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C169").Value = "Metropolitana"
$ws.Range("D169").Value = 44522
$ws.Range("E169").Value = 13
$ws.Range("F169").Value = "Fruta"
$ws.Range("G169").Value = 100107
$ws.Range("H169").Value = "Otros"
$ws.Range("I169").Value = 100107002
$ws.Range("J169").Value = "Chirimoya"
$ws.Range("K169").Value = "Cultivar IV Región"
$ws.Range("L169").Value = "Especial"
$ws.Range("M169").Value = 180
$ws.Range("N169").Value = 2200
$ws.Range("O169").Value = 2200
$ws.Range("P169").Value = 2200
$ws.Range("Q169").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R169").Value = "Provincia de Limarí"
$ws.Range("S169").Value = 2200
$ws.Range("T169").Value = 1

# Row 170: new weekly price entry (Extra (doble especial))
$ws.Range("A170").Value = 6
$ws.Range("B170").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C170").Value = "Metropolitana"
$ws.Range("D170").Value = 44522
$ws.Range("E170").Value = 13
$ws.Range("F170").Value = "Fruta"
$ws.Range("G170").Value = 100107
$ws.Range("H170").Value = "Otros"
$ws.Range("I170").Value = 100107002
$ws.Range("J170").Value = "Chirimoya"
$ws.Range("K170").Value = "Cultivar IV Región"
$ws.Range("L170").Value = "Extra (doble especial)"
$ws.Range("M170").Value = 170
$ws.Range("N170").Value = 2400
$ws.Range("O170").Value = 2400
$ws.Range("P170").Value = 2400
$ws.Range("Q170").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R170").Value = "Provincia de Limarí"
$ws.Range("S170").Value = 2400
$ws.Range("T170").Value = 1

# Row 171: new weekly price entry (Primera)
$ws.Range("A171").Value = 6
$ws.Range("B171").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C171").Value = "Metropolitana"
$ws.Range("D171").Value = 44522
$ws.Range("E171").Value = 13
$ws.Range("F171").Value = "Fruta"
$ws.Range("G171").Value = 100107
$ws.Range("H171").Value = "Otros"
$ws.Range("I171").Value = 100107002
$ws.Range("J171").Value = "Chirimoya"
$ws.Range("K171").Value = "Cultivar IV Región"
$ws.Range("L171").Value = "Primera"
$ws.Range("M171").Value = 200
$ws.Range("N171").Value = 1900
$ws.Range("O171").Value = 1900
$ws.Range("P171").Value = 1900
$ws.Range("Q171").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R171").Value = "Provincia de Limarí"
$ws.Range("S171").Value = 1900
$ws.Range("T171").Value = 1
